$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.791.47'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.92%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.479.35'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.11%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.24%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.513'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.36%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.144'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.18%  '

$ws.Range("E10").Value = '  -1.49%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '4.96'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.335'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.53%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.933.91'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.09%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.32'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.80%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '67.687.70'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.98%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000170'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.493.64'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.90%  '

$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.84'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.55%  '

$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.78%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '348.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.84%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.54%  '

$ws.Range("E22").Value = '  -0.18%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.78'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.94%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.22%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -7.09%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.82'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.55%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.587.21'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.60%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.995'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.35%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0891'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.64%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '496.83'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.71'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.19%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.25'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.25%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.76'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.53%  '

$ws.Range("E34").Value = '  +0.06%  '

$ws.Range("B35").Value = 'Monero'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '163.07'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.44%  '

$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.120'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.48%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.63'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.39%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.29'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.78%  '

$ws.Range("E39").Value = '  -0.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.29'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.72'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.28%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.326'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.96%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.79'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.41'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.23%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '147.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.89%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.52'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.05%  '

$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.510'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.00%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0253'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.54%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0736'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.38%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.56'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.73%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.576'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.39%  '
